$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26: new "day 25" entry ---
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 43855
$ws.Range("B26").NumberFormat = $ws.Range("B25").NumberFormat

$ws.Range("C26").Value = "Surah Baqarah, 258 - 264"
$ws.Range("F26").Value = "Sacrifice of Ibrahim (as), Foundation of Kaaba, Love for Allah (swt), Life of Ibrahim (as)"

$ws.Range("D26").Value = @'
h1: What is my worth?
p: Tall men. Really tall lived on this planet long ago in history. They used to live 1000 years. Some generations were really advanced. They could find the location of angels looking at stars. Some built pyramids. Some build planes that fly. Some went to moon. A variety of men have lived and died. Some achieved great successes and secured huge lands. Ruled the planet. Some died a failure. Among all these generations, who is it that really mattered to the Creator of this planet?
p.b-left: People who prostrated to Allah (swt) the most?
p.b-left: People who ruled with justice and eliminated injustices. Did they matter the most?
p.b-left: People who taught Quran and spread its word? Did they matter most? 
h3: Who was it that mattered?
p: <b>Ibrahim (as)</b>. 
p.note: I am going to try to picture the situation he faced, using today’s situation. Try to be with me.
p: The present day, we have a good normal life. We have secured a financial stability and our parents, wife and kids are living a fine life. Kids will get good education, parents will get good treatment in hospitals and wife will get good freedom. Everything is normal in any normal person’s life. 
p: Making it abnormal would be a chaos. It will be a down hill run and no coming back. <b>Raising voice for injustice happening somewhere out in the courts</b> would be abnormal. Trying to feed a kid for lifetime while<b> compromising my own kids studies</b> would be abnormal. <b>Selling all of our savings out and giving it in sadaqah</b> is abnormal. <b>Fighting to protect the underprivileged</b> in a normal life would be abnormal. 
p: Running down this path of abnormality is not recommended, not allowed. Dying in the path of Allah is not recommended by our relatives and friends. Shahadat (Martyrdoom) is considered a loss to a young wife. Praying for it, or saying to achieve it out loud, is bad.
p: <b>Lets say, someone raises his voice for justice in an unjust society. What will happen? </b>
p: He will be kicked out the 1st day from his home. He will run around for protection out in the open. No one will give him protection if he keeps acting “abnormal”. 
p: So he will break down and fall, may be tomorrow or some days ahead. He will surrender to the norms of society. He will try to act a little normal that he is given a shelter and survive.
p: If he sticks to “<b>raising voice against injustice</b>” and “<b>does not settle for any less he determined to</b>”, he will be attacked harder by people in power. His own relatives first, followed by colleagues and neighbours will try to pull him down out of “love for his wellbeing”. 
p: <b>Still if he does not give up. </b>
p: He will be given electric shocks. He will be put behind bars. He will be made bedridden forcefully. He will be dented to fit in the definition of a<b> ‘Normal Life’</b>.
h3: Who was Ibrahim (as)?
p: Ibrahim (as) was somebody who after the electric shocks, the phase behind bars and the hard hitting; did not settle. He kept on running towards the fire, was thrown into it. He was made to question his own sanity. Just a dream, made him put his own son down for “straight to the point”, ”no questions asked” sacrifice. 
p: From being the only 1 on this planet to believe in the might of Allah (swt), he made it to 1.8 Billion people believing in the might of Creator. We face the stone, he moved. We face the city, he built. We love the sacrifice, he made. What to talk about me or you or any other muslim. Our Prophet went into isolation to find Allah, Ibrahim found. He looked towards Allah (swt) to make the Kaaba, a Kibla, Ibrahim built. 
h3: How hard did Ibrahim (as) fought?
p: Some beaten, declared abnormal goes to Chief Justice of some state and questions his credibility. He will be humiliated, threatened and pushed further lower. Ibrahim (as) did not go to Chief Justice. He went straight to the ruler ‘Namrood’. Ibrahim (as) had neither power nor reputation to settle an argument with him. He straight on went into questioning the credibility of the ruler. 
quote: Have you not considered the one who argued with Abraham about his Lord [merely] because Allah had given him kingship? When Abraham said, "My Lord is the one who gives life and causes death," he said, "I give life and cause death." Abraham said, "Indeed, Allah brings up the sun from the east, so bring it up from the west." So the disbeliever was overwhelmed [by astonishment], and Allah does not guide the wrongdoing people. <br>- Surah Baqarah verse 258
p: More sufferings his way. Did he give up? Did he not fall into fire for Allah (swt)? Did he not lay his own son a sacrifice? Did he not settle for nothing but truth?
p: Somebody did it for the love of Allah (swt). Where do I stand? After saying 5 prayers and writing articles and talking high. Still at 0. My worth? Like anybody’s worth, is 0 when it comes to comparing my superficial love with Ibrahim (as)’s love.
p: <b>May Allah (swt) accept our efforts and make us brave enough that we follow the path of Ibrahim (as). May we make our Creator proud of us.. Amen </b>
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
'@

$ws.Range("E26").Value = "Qasim Ali"

$ws.Rows.Item(26).RowHeight = 409.6

# --- Update active sheet view: scroll to row 26, select D26 ---
$ws.Range("D26").Select() | Out-Null
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 1

Write-Host "Row 26 populated."
